# Updates the "base de datos" table (rows 16-23) on Hoja1:
#  - Removes the previous account-statement periods for
#    JORGE ENRIQUE GARCIA ORTEGA (2308-2402) and JAIRO MERCADO BUSTAMANTE (2308)
#  - Adds new rows for JAIRO MERCADO BUSTAMANTE (2308, 2402) and
#    JORGE ENRIQUE GARCIA ORTEGA (2401, 2312, 2311, 2310, 2309, 2308)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rows = @(
    @{ Row = 16; DocType = "CC"; DocNum = "73203395";  Name = "JAIRO MERCADO BUSTAMANTE";    Period = "2308"; ValorMora = 46400; Salario = 737717  },
    @{ Row = 17; DocType = "CC"; DocNum = "73203395";  Name = "JAIRO MERCADO BUSTAMANTE";    Period = "2402"; ValorMora = 24000; Salario = 1160000 },
    @{ Row = 18; DocType = "CC"; DocNum = "1102816018"; Name = "JORGE ENRIQUE GARCIA ORTEGA"; Period = "2401"; ValorMora = 40000; Salario = 1160000 },
    @{ Row = 19; DocType = "CC"; DocNum = "1102816018"; Name = "JORGE ENRIQUE GARCIA ORTEGA"; Period = "2312"; ValorMora = 46400; Salario = 1160000 },
    @{ Row = 20; DocType = "CC"; DocNum = "1102816018"; Name = "JORGE ENRIQUE GARCIA ORTEGA"; Period = "2311"; ValorMora = 46400; Salario = 1160000 },
    @{ Row = 21; DocType = "CC"; DocNum = "1102816018"; Name = "JORGE ENRIQUE GARCIA ORTEGA"; Period = "2310"; ValorMora = 46400; Salario = 1160000 },
    @{ Row = 22; DocType = "CC"; DocNum = "1102816018"; Name = "JORGE ENRIQUE GARCIA ORTEGA"; Period = "2309"; ValorMora = 46400; Salario = 1160000 },
    @{ Row = 23; DocType = "CC"; DocNum = "1102816018"; Name = "JORGE ENRIQUE GARCIA ORTEGA"; Period = "2308"; ValorMora = 46400; Salario = 1160000 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.DocType
    $ws.Cells.Item($r.Row, 3).Value = $r.DocNum
    $ws.Cells.Item($r.Row, 4).Value = $r.Name
    $ws.Cells.Item($r.Row, 5).Value = $r.Period
    $ws.Cells.Item($r.Row, 6).Value = $r.ValorMora
    $ws.Cells.Item($r.Row, 7).Value = $r.Salario
}
